$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column K (최종점수 / final score) values
$ws.Range("K2").Value = 62.1
$ws.Range("K3").Value = 55.7

# Update column N (MACRO_SCORE) values
$ws.Range("N2").Value = 85.82376350509293
$ws.Range("N3").Value = 85.82376350509293
